$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "10+80="
$t.Cell(1, 2).Range.Text = "88-67="
$t.Cell(1, 3).Range.Text = "9+47="
$t.Cell(1, 4).Range.Text = "2+40="
$t.Cell(1, 5).Range.Text = "23+52="
$t.Cell(2, 1).Range.Text = "88-39="
$t.Cell(2, 2).Range.Text = "38+26="
$t.Cell(2, 3).Range.Text = "49+1="
$t.Cell(2, 4).Range.Text = "53-2="
$t.Cell(2, 5).Range.Text = "15+51="
$t.Cell(3, 1).Range.Text = "22+4="
$t.Cell(3, 2).Range.Text = "21+47="
$t.Cell(3, 3).Range.Text = "92-32="
$t.Cell(3, 4).Range.Text = "43+33="
$t.Cell(3, 5).Range.Text = "47-9="
$t.Cell(4, 1).Range.Text = "44-2="
$t.Cell(4, 2).Range.Text = "44+28="
$t.Cell(4, 3).Range.Text = "30-13="
$t.Cell(4, 4).Range.Text = "9+13="
$t.Cell(4, 5).Range.Text = "42-1="
$t.Cell(5, 1).Range.Text = "56-49="
$t.Cell(5, 2).Range.Text = "10+71="
$t.Cell(5, 3).Range.Text = "58-38="
$t.Cell(5, 4).Range.Text = "73+17="
$t.Cell(5, 5).Range.Text = "89-61="
$t.Cell(6, 1).Range.Text = "4+61="
$t.Cell(6, 2).Range.Text = "35-10="
$t.Cell(6, 3).Range.Text = "47+9="
$t.Cell(6, 4).Range.Text = "38+19="
$t.Cell(6, 5).Range.Text = "72+11="
$t.Cell(7, 1).Range.Text = "83-50="
$t.Cell(7, 2).Range.Text = "59+39="
$t.Cell(7, 3).Range.Text = "86-45="
$t.Cell(7, 4).Range.Text = "7-2="
$t.Cell(7, 5).Range.Text = "76-15="
$t.Cell(8, 1).Range.Text = "25+65="
$t.Cell(8, 2).Range.Text = "36+55="
$t.Cell(8, 3).Range.Text = "48-4="
$t.Cell(8, 4).Range.Text = "66-16="
$t.Cell(8, 5).Range.Text = "20-5="
$t.Cell(9, 1).Range.Text = "91-20="
$t.Cell(9, 2).Range.Text = "59+22="
$t.Cell(9, 3).Range.Text = "91-33="
$t.Cell(9, 4).Range.Text = "22-15="
$t.Cell(9, 5).Range.Text = "22+52="
$t.Cell(10, 1).Range.Text = "72-61="
$t.Cell(10, 2).Range.Text = "35+26="
$t.Cell(10, 3).Range.Text = "7-2="
$t.Cell(10, 4).Range.Text = "49-26="
$t.Cell(10, 5).Range.Text = "36-35="
$t.Cell(11, 1).Range.Text = "20+7="
$t.Cell(11, 2).Range.Text = "39+55="
$t.Cell(11, 3).Range.Text = "0+47="
$t.Cell(11, 4).Range.Text = "79-54="
$t.Cell(11, 5).Range.Text = "69-41="
$t.Cell(12, 1).Range.Text = "77-62="
$t.Cell(12, 2).Range.Text = "57-54="
$t.Cell(12, 3).Range.Text = "26-23="
$t.Cell(12, 4).Range.Text = "10+89="
$t.Cell(12, 5).Range.Text = "89-33="
$t.Cell(13, 1).Range.Text = "12+49="
$t.Cell(13, 2).Range.Text = "89-88="
$t.Cell(13, 3).Range.Text = "57-29="
$t.Cell(13, 4).Range.Text = "60+24="
$t.Cell(13, 5).Range.Text = "73-45="
$t.Cell(14, 1).Range.Text = "77-44="
$t.Cell(14, 2).Range.Text = "94-59="
$t.Cell(14, 3).Range.Text = "15+50="
$t.Cell(14, 4).Range.Text = "38+38="
$t.Cell(14, 5).Range.Text = "3+68="
$t.Cell(15, 1).Range.Text = "25+6="
$t.Cell(15, 2).Range.Text = "53+3="
$t.Cell(15, 3).Range.Text = "76-35="
$t.Cell(15, 4).Range.Text = "10+25="
$t.Cell(15, 5).Range.Text = "92-43="
$t.Cell(16, 1).Range.Text = "55-31="
$t.Cell(16, 2).Range.Text = "15-9="
$t.Cell(16, 3).Range.Text = "1+22="
$t.Cell(16, 4).Range.Text = "37+4="
$t.Cell(16, 5).Range.Text = "89-58="
$t.Cell(17, 1).Range.Text = "68-6="
$t.Cell(17, 2).Range.Text = "12+41="
$t.Cell(17, 3).Range.Text = "30+0="
$t.Cell(17, 4).Range.Text = "96-60="
$t.Cell(17, 5).Range.Text = "48-1="
$t.Cell(18, 1).Range.Text = "10+44="
$t.Cell(18, 2).Range.Text = "94-41="
$t.Cell(18, 3).Range.Text = "5+39="
$t.Cell(18, 4).Range.Text = "82+2="
$t.Cell(18, 5).Range.Text = "92-13="
$t.Cell(19, 1).Range.Text = "89-65="
$t.Cell(19, 2).Range.Text = "40-14="
$t.Cell(19, 3).Range.Text = "28+4="
$t.Cell(19, 4).Range.Text = "38+27="
$t.Cell(19, 5).Range.Text = "2+85="
$t.Cell(20, 1).Range.Text = "65-3="
$t.Cell(20, 2).Range.Text = "31-2="
$t.Cell(20, 3).Range.Text = "79-8="
$t.Cell(20, 4).Range.Text = "63-34="
$t.Cell(20, 5).Range.Text = "66-24="
